# Update Success/Prediction/Error (and Cross Entropy Loss / Success %) values
# for the NODE-CNN classification results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ C = 0; D = 0.920025824800212;  E = 0.920025824800212 }
    3  = @{ C = 0; D = 0.5023501437370876; E = 0.5023501437370876 }
    4  = @{ C = 0; D = 0.8735947769639344; E = 0.8735947769639344 }
    5  = @{ C = 1; D = 0.2900510517540365; E = 0.2900510517540365 }
    6  = @{ C = 1; D = 0.2509929039701788; E = 0.2509929039701788 }
    7  = @{ C = 0; D = 0.3791828664405087; E = 0.6208171335594912 }
    8  = @{ C = 0; D = 0.02478555272305007; E = 0.9752144472769499 }
    9  = @{ C = 0; D = 0.2959387053810721; E = 0.7040612946189279 }
    10 = @{ C = 1; D = 0.5896979176579943; E = 0.4103020823420057 }
    11 = @{ C = 0; D = 0.4936121180219135; E = 0.5063878819780865; F = 1.304272532463074; G = 0.3 }
    12 = @{ C = 0; D = 0.956874411011657;  E = 0.956874411011657 }
    13 = @{ C = 0; D = 0.6219649011393844; E = 0.6219649011393844 }
    14 = @{ C = 0; D = 0.9015223530272547; E = 0.9015223530272547 }
    15 = @{ C = 0; D = 0.5346766011167032; E = 0.5346766011167032 }
    16 = @{ C = 1; D = 0.230688658512838;  E = 0.230688658512838 }
    17 = @{ C = 1; D = 0.512905428161204;  E = 0.487094571838796 }
    18 = @{ C = 0; D = 0.004239941232161234; E = 0.9957600587678388 }
    19 = @{ C = 0; D = 0.3322931244332533; E = 0.6677068755667467 }
    20 = @{ C = 1; D = 0.7414485711831754; E = 0.2585514288168246 }
    21 = @{ C = 0; D = 0.6246586447050348; E = 0.3753413552949652; F = 1.546392202377319; G = 0.4 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = [bool]$vals.C
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
    if ($vals.ContainsKey('F')) {
        $ws.Cells.Item($row, 6).Value = $vals.F
    }
    if ($vals.ContainsKey('G')) {
        $ws.Cells.Item($row, 7).Value = $vals.G
    }
}
